# Updates the "cryptos" price table (Sun Jul 21 15:48:08 UTC 2024 refresh).
#
# Columns: B=Coin, C=Link, D=Price, E=Volume(1h). D/E are stored as literal
# text (the sheet has no numeric formatting for them, and values such as
# "67.360.19" / "4.117.24" aren't valid numbers anyway - they're thousands-
# dot-grouped price strings). Whenever a new Price value DOES look like a
# genuine number (e.g. "597.53", "8.00", "0.180"), Excel's COM Value setter
# would silently reinterpret it as a numeric cell and normalise away
# significant trailing/grouping digits (e.g. "8.00" -> 8, "0.180" -> 0.18).
# To keep these as plain text - matching the original inline-string cells -
# the cell is switched to Text number format ("@") before the value is
# written. Volume(1h) values always carry two leading/trailing spaces and a
# trailing "%" so Excel never parses them as numbers and no such guard is
# needed there.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '67.193.25'
$ws.Range("E2").Value = '  +0.74%  '

$ws.Range("D3").Value = '3.503.50'
$ws.Range("E3").Value = '  +0.08%  '

$ws.Range("E4").Value = '  +0.08%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '597.53'
$ws.Range("E5").Value = '  +0.65%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '175.56'
$ws.Range("E6").Value = '  +3.83%  '

$ws.Range("E7").Value = '  +0.04%  '

$ws.Range("E8").Value = '  -0.58%  '

$ws.Range("E9").Value = '  -0.89%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '7.17'
$ws.Range("E10").Value = '  -2.06%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.431'
$ws.Range("E11").Value = '  -0.54%  '

$ws.Range("D12").Value = '4.115.56'
$ws.Range("E12").Value = '  +0.24%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '30.73'
$ws.Range("E13").Value = '  +8.46%  '

$ws.Range("E14").Value = '  +0.18%  '

$ws.Range("D15").Value = '67.255.88'
$ws.Range("E15").Value = '  +0.81%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.0000179'
$ws.Range("E16").Value = '  -1.41%  '

$ws.Range("D17").Value = '3.512.47'
$ws.Range("E17").Value = '  +0.30%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '6.30'
$ws.Range("E18").Value = '  -0.47%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '14.48'
$ws.Range("E19").Value = '  +3.16%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '393.96'
$ws.Range("E20").Value = '  -0.57%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '8.00'
$ws.Range("E21").Value = '  +0.17%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '73.45'
$ws.Range("E22").Value = '  -0.02%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.999'
$ws.Range("E23").Value = '  +0.15%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.537'
$ws.Range("E24").Value = '  +0.16%  '

$ws.Range("E25").Value = '  -0.57%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.0000122'
$ws.Range("E26").Value = '  -0.07%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '10.19'
$ws.Range("E27").Value = '  +0.12%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.180'
$ws.Range("E28").Value = '  -0.14%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.996'
$ws.Range("E29").Value = '  -0.34%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '6.14'
$ws.Range("E30").Value = '  -2.63%  '

$ws.Range("E31").Value = '  -2.32%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '2.06'
$ws.Range("E32").Value = '  -0.38%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '23.65'
$ws.Range("E33").Value = '  -0.66%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '7.37'
$ws.Range("E34").Value = '  -0.26%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.64'
$ws.Range("E35").Value = '  +1.26%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '163.41'
$ws.Range("E36").Value = '  +0.49%  '

$ws.Range("E37").Value = '  -1.85%  '

$ws.Range("E38").Value = '  +0.61%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '7.01'
$ws.Range("E39").Value = '  +2.70%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '27.64'
$ws.Range("E40").Value = '  +2.03%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '4.66'
$ws.Range("E41").Value = '  -0.45%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.0732'
$ws.Range("E42").Value = '  -1.72%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '26.11'
$ws.Range("E43").Value = '  -1.61%  '

$ws.Range("D44").Value = '2.798.92'
$ws.Range("E44").Value = '  +0.04%  '

$ws.Range("E45").Value = '  -0.86%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '42.47'
$ws.Range("E46").Value = '  -0.98%  '

$ws.Range("E47").Value = '  -2.95%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '340.93'
$ws.Range("E48").Value = '  -0.38%  '

$ws.Range("E49").Value = '  -1.01%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '33.65'
$ws.Range("E50").Value = '  -0.94%  '

$ws.Range("B51").Value = 'Cosmos'
$ws.Range("C51").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '6.44'
$ws.Range("E51").Value = '  -1.03%  '

